# Auto-generated edit script: updates cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.806.30"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "1.599.73"
$ws.Range("E3").Value = "  -2.25%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'208.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.52%  "
$ws.Range("D6").Value = "'1.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.480"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.71%  "
$ws.Range("E8").Value = "  -2.10%  "
$ws.Range("D9").Value = "'0.0610"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("D10").Value = "'17.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.15%  "
$ws.Range("D12").Value = "1.824.19"
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("D13").Value = "1.603.71"
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("E14").Value = "  -3.68%  "
$ws.Range("E15").Value = "  -3.97%  "
$ws.Range("D16").Value = "25.809.48"
$ws.Range("D17").Value = "'60.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("D18").Value = "0.0₃0717"
$ws.Range("E18").Value = "  -3.51%  "
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").Value = "'189.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("D21").Value = "'4.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.42%  "
$ws.Range("E22").Value = "  -3.89%  "
$ws.Range("D23").Value = "'5.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D24").Value = "'0.129"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.21%  "
$ws.Range("B25").Value = "BinanceUSD"
$ws.Range("C25").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D25").Value = "'1.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("D26").Value = "'140.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.97%  "
$ws.Range("E27").Value = "  -4.53%  "
$ws.Range("D28").Value = "'6.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.86%  "
$ws.Range("E29").Value = "  -1.79%  "
$ws.Range("E30").Value = "  -4.11%  "
$ws.Range("E31").Value = "  -3.22%  "
$ws.Range("E32").Value = "  -2.70%  "
$ws.Range("D33").Value = "'3.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.82%  "
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("E35").Value = "  -1.43%  "
$ws.Range("D36").Value = "1.095.55"
$ws.Range("E36").Value = "  -3.55%  "
$ws.Range("D37").Value = "'2.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("B38").Value = "PaxDollar"
$ws.Range("C38").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("E39").Value = "  -2.31%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'0.792"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.57%  "
$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D41").Value = "'0.497"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.36%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'95.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.86%  "
$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").Value = "1.735.99"
$ws.Range("E43").Value = "  -2.10%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.24%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'0.743"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.47%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0112"
$ws.Range("E46").Value = "  -2.31%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'53.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.51%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.0512"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.00%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.43"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.99%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.410"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.25%  "
